# Updates cryptos list data (prices / 1h volume changes) to match the
# latest scrape, including a reordering of the WrappedBTC /
# WrappedliquidstakedEther2.0 rows (15 and 16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.503.74'
$ws.Range("E2").Value = '  +3.25%  '
$ws.Range("D3").Value = '2.333.20'
$ws.Range("E3").Value = '  +1.37%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '545.57'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.37%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '131.47'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("E7").Value = '  +0.01%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.578'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.34%  '
$ws.Range("D9").Value = '2.330.26'
$ws.Range("E9").Value = '  +1.23%  '
$ws.Range("E10").Value = '  +0.83%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '5.54'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("E13").Value = '  +0.46%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '23.68'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.747.39'
$ws.Range("E15").Value = '  +1.32%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '60.457.28'
$ws.Range("E16").Value = '  +3.36%  '
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '2.332.93'
$ws.Range("E18").Value = '  +1.32%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '10.61'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.31%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.15'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.38%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '315.38'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.26%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.67'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.91%  '
$ws.Range("E23").Value = '  -0.04%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '64.00'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.46%  '
$ws.Range("E25").Value = '  +2.15%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.03%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.87'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.08%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.36'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +5.12%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.20'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +10.19%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '173.29'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("E31").Value = '  +1.21%  '
$ws.Range("D32").Value = '0.0₃0735'
$ws.Range("E32").Value = '  +1.35%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.95'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.31%  '
$ws.Range("E34").Value = '  +11.31%  '
$ws.Range("E35").Value = '  -0.89%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '17.85'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("E38").Value = '  +0.10%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '4.08'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.21%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '322.26'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +10.86%  '
$ws.Range("E41").Value = '  +1.99%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '38.01'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.82%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '137.89'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.22%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '3.50'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.19%  '
$ws.Range("E45").Value = '  -1.21%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '19.24'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +5.37%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.565'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.61%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0497'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.26%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.0214'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.75%  '
$ws.Range("D50").Value = '0.0₆0213'
$ws.Range("E50").Value = '  +15.36%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '11.03'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.72%  '
